$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving the cell's
# original (default) style - COM auto-converts numeric-looking
# strings to actual numbers unless the cell is explicitly text-
# formatted first; we then restore the style so no stray formatting
# change is introduced.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "54.518.73"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.289.87"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "503.65"
$ws.Range("E5").Value = "  +2.04%  "
Set-TextValue "D6" "130.72"
$ws.Range("E6").Value = "  +2.76%  "
Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  -0.22%  "
Set-TextValue "D8" "0.530"
$ws.Range("E8").Value = "  +0.53%  "
Set-TextValue "D9" "0.0958"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("E10").Value = "  +0.56%  "
Set-TextValue "D11" "0.340"
$ws.Range("E11").Value = "  +5.05%  "
Set-TextValue "D12" "4.75"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").Value = "2.700.03"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +6.33%  "
$ws.Range("D15").Value = "54.482.83"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "2.305.49"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("E18").Value = "  +3.17%  "
Set-TextValue "D19" "4.17"
$ws.Range("E19").Value = "  +2.95%  "
Set-TextValue "D20" "304.93"
$ws.Range("E20").Value = "  +0.55%  "
Set-TextValue "D21" "6.36"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("E22").Value = "  +0.02%  "
Set-TextValue "D23" "61.99"
$ws.Range("E23").Value = "  -2.49%  "
Set-TextValue "D24" "0.997"
$ws.Range("E24").Value = "  -0.43%  "
Set-TextValue "D25" "0.151"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  +3.48%  "
Set-TextValue "D27" "171.18"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").Value = "0.0₃0697"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  -0.03%  "
Set-TextValue "D33" "17.84"
$ws.Range("E33").Value = "  +1.62%  "
Set-TextValue "D34" "0.967"
$ws.Range("E34").Value = "  +11.03%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +3.33%  "
Set-TextValue "D38" "0.377"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D40" "5.08"
$ws.Range("E40").Value = "  +5.83%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D41" "3.40"
$ws.Range("E41").Value = "  +1.81%  "
Set-TextValue "D42" "126.52"
$ws.Range("E42").Value = "  -1.60%  "
Set-TextValue "D43" "0.0497"
$ws.Range("E43").Value = "  +3.86%  "
Set-TextValue "D44" "0.0899"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  +1.31%  "
Set-TextValue "D46" "242.69"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").Value = "  +0.76%  "
Set-TextValue "D50" "16.48"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("E51").Value = "  +1.99%  "
